# Adding extract values from excel file to add book api example
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("testdata")

# New row of test data appended after the existing 5 rows.
$ws.Range("A6").Value = "RestAddbook"
$ws.Range("B6").Value = "RestAssured"
$ws.Range("C6").Value = "asfasf"
$ws.Range("D6").Value = 3221

# Move selection/active cell to match the saved view state.
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 4
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("H6").Select()
